# AFDP-2522: Update of foia drools rule files
# The FOIA "Billing"/"Hold" queue enter-date rules had their CONDITION
# expressions simplified:
#   - the defensive `?.` safe-navigation on `queue?.name` was dropped
#     (now `queue.name`)
#   - the "nullify" rules' redundant second clause
#     (`&& billingEnterDate != null` / `&& holdEnterDate != null`) was
#     removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28 - Set Billing Enter Date
$ws.Range("C28").Value = "queue.name == 'Billing' && billingEnterDate == null"

# Row 29 - Nullify Billing Enter Date
$ws.Range("C29").Value = "queue.name != 'Billing'"

# Row 30 - Set Hold Enter Date
$ws.Range("C30").Value = "queue.name == 'Hold' && holdEnterDate == null"

# Row 31 - Nullify Hold Enter Date
$ws.Range("C31").Value = "queue.name != 'Hold'"

# Move the active selection to D31 to match the saved cursor position.
[void]$ws.Range("D31").Select()
